$wb = $excel.ActiveWorkbook

# --- Sheet: Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

# pre-format the text columns so assigning numeric-looking strings keeps them as text
$ws.Range("B2:B4").NumberFormat = "@"
$ws.Range("D2:F4").NumberFormat = "@"

# Row 2 (J_0_L0_v) -- set in left-to-right column order: A, B, (C untouched), D, E, F
$ws.Range("A2").Value = "-2.5315837802100924 - 2x_1 + 1.8678224246431105y_1 - 0.14263326722005942y_2"
$ws.Range("B2").Value = "5.0315837802100924"
$ws.Range("D2").Value = "0.96"
$ws.Range("E2").Value = "3.5"
$ws.Range("F2").Value = "4.3"

# Row 3 (J_0_LP_v)
$ws.Range("A3").Value = "10.204623175446986 + x_1 - 3x_2 - 1.7499928618506757y_1 - 0.7289086105030773y_2"
$ws.Range("B3").Value = "-12.204623175446986"
$ws.Range("D3").Value = "0.9"
$ws.Range("E3").Value = "8.5"
$ws.Range("F3").Value = "5.699999999999999"

# Row 4 (J_Ne_L0_v)
$ws.Range("A4").Value = "-15.949985308161315 + x_1 + x_2 + 1.7884864957516098y_1 + 1.7669384656823133y_2"
$ws.Range("B4").Value = "13.199985308161315"
$ws.Range("D4").Value = "0.28"
$ws.Range("E4").Value = "7.6"
$ws.Range("F4").Value = "8.8"

# --- Sheet: Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2:D2").NumberFormat = "@"
$ws.Range("A2").Value = "2.55"
$ws.Range("B2").Value = "1.25"
$ws.Range("C2").Value = "4.25"
$ws.Range("D2").Value = "2.15"

# --- Sheet: Vector_bf (note: "Vector_bf" and "Vector_BF" differ only by case,
#     so Worksheets.Item(name) is ambiguous; use positional index instead) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2:A3").NumberFormat = "@"
$ws.Range("A2").Value = "3.2811078291977713"
$ws.Range("A3").Value = "-0.7017970844070212"

# --- Sheet: Vector_BF ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2:A5").NumberFormat = "@"
$ws.Range("A2").Value = "-7.1"
$ws.Range("A3").Value = "16.9"
$ws.Range("A4").Value = "-5.754936528232378"
$ws.Range("A5").Value = "-6.733792714639216"

# --- Sheet: Vector_Alpha ---
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 2.4899999999999998
$ws.Range("A3").Value = 2.46
